$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark that currently sits on the
#    "Participantes" heading paragraph.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2) Insert a new list item "Alba Vallés Esteban" (bold) right after
#    "Francisco Javier Nogueras Iso" and before "Enrique Ruiz Flores",
#    re-using the same list/paragraph formatting, and re-attach the
#    "_GoBack" bookmark as a zero-length bookmark right after the new
#    run's text (before the paragraph mark).
# ------------------------------------------------------------------
$franciscoPara = $d.Paragraphs(3)
$franciscoPara.Range.InsertParagraphAfter()

$albaPara = $d.Paragraphs(4)
$albaPara.Range.Text = "Alba Vallés EstebanZZ"

$albaRange = $albaPara.Range
$bmPos = $albaRange.End - 3
$safeBmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $safeBmRange)

$dummyRange = $d.Range($bmPos, $bmPos + 2)
$dummyRange.Delete()

# ------------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from the run that starts
#    "Se reconsiderará la " to the run that starts
#    "Se reescribirá el RNF-5 para ".
# ------------------------------------------------------------------
$rnf5Para = $null
$reconsideraraPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.Contains("RNF-5")) {
        $rnf5Para = $p
    }
    if ($t.Contains("pasarela de pagos")) {
        $reconsideraraPara = $p
    }
}

$rnf5Range = $rnf5Para.Range
$rnf5Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="031C11A4" w14:textId="05626993" w:rsidR="000F1FEC" w:rsidRDefault="00853D59" w:rsidP="000A3FBC"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Se reescribirá el RNF-5 para </w:t></w:r><w:r w:rsidR="004C2EB3"><w:t>que sea más concreto</w:t></w:r><w:r w:rsidR="00833936"><w:t>.</w:t></w:r></w:p>'
$rnf5Range.InsertXML($rnf5Xml)

$reconsideraraRange = $reconsideraraPara.Range
$reconsideraraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2B78E1AE" w14:textId="124C8FF0" w:rsidR="00833936" w:rsidRDefault="000B63CA" w:rsidP="000A3FBC"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Se reconsiderará la </w:t></w:r><w:r w:rsidR="00105B10"><w:t>pasarela de pagos, puesto que puede resultar complicada</w:t></w:r><w:r w:rsidR="006D07D7"><w:t>.</w:t></w:r></w:p>'
$reconsideraraRange.InsertXML($reconsideraraXml)

Write-Output "edit complete"
